$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 336, shifting existing rows 336-410
# down to 339-413 (preserving their values/formatting automatically).
$ws.Rows("336:338").Insert()

# --- Populate new row 336 ---
$ws.Range("A336").Value = 5
$ws.Range("B336").Value = "Macroferia Regional de Talca"
$ws.Range("C336").Value = "Maule"
$ws.Range("D336").Value = 44551
$ws.Range("E336").Value = 7
$ws.Range("F336").Value = "Fruta"
$ws.Range("G336").Value = 100101
$ws.Range("H336").Value = "Berries"
$ws.Range("I336").Value = 100112025
$ws.Range("J336").Value = "Frutilla"
$ws.Range("K336").Value = "Sin especificar"
$ws.Range("L336").Value = "Especial"
$ws.Range("M336").Value = 150
$ws.Range("N336").Value = 8000
$ws.Range("O336").Value = 8000
$ws.Range("P336").Value = 8000
$ws.Range("Q336").Value = "$/bandeja 7 kilos"
$ws.Range("R336").Value = "Provincia de Melipilla"
$ws.Range("S336").Value = 1143
$ws.Range("T336").Value = 7

# --- Populate new row 337 ---
$ws.Range("A337").Value = 5
$ws.Range("B337").Value = "Macroferia Regional de Talca"
$ws.Range("C337").Value = "Maule"
$ws.Range("D337").Value = 44551
$ws.Range("E337").Value = 7
$ws.Range("F337").Value = "Fruta"
$ws.Range("G337").Value = 100101
$ws.Range("H337").Value = "Berries"
$ws.Range("I337").Value = 100112025
$ws.Range("J337").Value = "Frutilla"
$ws.Range("K337").Value = "Sin especificar"
$ws.Range("L337").Value = "Especial"
$ws.Range("M337").Value = 150
$ws.Range("N337").Value = 7000
$ws.Range("O337").Value = 7000
$ws.Range("P337").Value = 7000
$ws.Range("Q337").Value = "$/bandeja 7 kilos"
$ws.Range("R337").Value = "Región del Maule"
$ws.Range("S337").Value = 1000
$ws.Range("T337").Value = 7

# --- Populate new row 338 ---
$ws.Range("A338").Value = 5
$ws.Range("B338").Value = "Macroferia Regional de Talca"
$ws.Range("C338").Value = "Maule"
$ws.Range("D338").Value = 44551
$ws.Range("E338").Value = 7
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100101
$ws.Range("H338").Value = "Berries"
$ws.Range("I338").Value = 100112025
$ws.Range("J338").Value = "Frutilla"
$ws.Range("K338").Value = "Sin especificar"
$ws.Range("L338").Value = "Especial"
$ws.Range("M338").Value = 150
$ws.Range("N338").Value = 8000
$ws.Range("O338").Value = 8000
$ws.Range("P338").Value = 8000
$ws.Range("Q338").Value = "$/caja 7 kilos"
$ws.Range("R338").Value = "Región del Maule"
$ws.Range("S338").Value = 1143
$ws.Range("T338").Value = 7

Write-Output "Done"
